# Electricity Dispatch Logit Exponent update:
# Align EDLE value with 2020 generation / BGDPbES guaranteed-dispatch change.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsEDLE  = $wb.Worksheets.Item("EDLE")

# Update the logit exponent value (B2) on the EDLE sheet.
$wsEDLE.Range("B2").Value = -3

# Reflect the cell selection left on the EDLE sheet (B3) after the edit,
# then restore the original active sheet (About) so the workbook still
# opens on the same tab it did before.
$wsEDLE.Activate()
$wsEDLE.Range("B3").Select()
$wsAbout.Activate()

$wb.Save()
